$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells (coin name / link) - safe to set directly ---
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'

# --- Numeric-looking text cells (price / volume %) must be forced to Text ---
# so Excel does not auto-convert them to real numbers (which would also drop
# significant trailing zeros, e.g. "8.730" -> 8.73, or turn "-1.30%" into a
# real percentage number instead of the literal string "-1.30%").
$forceCells = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9",
    "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "E16", "D17",
    "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "E24",
    "D25", "E25", "D26", "E26", "D27", "E27", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42",
    "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49",
    "D50", "E50", "D51", "E51"
)
foreach ($addr in $forceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '327.39'
$ws.Range("E2").Value = '-1.30%'
$ws.Range("D3").Value = '43.56'
$ws.Range("E3").Value = '5.02%'
$ws.Range("D4").Value = '5.483'
$ws.Range("E4").Value = '-3.08%'
$ws.Range("D5").Value = '0.08134'
$ws.Range("E5").Value = '-2.48%'
$ws.Range("D6").Value = '8.713'
$ws.Range("E6").Value = '-0.88%'
$ws.Range("D7").Value = '4.324'
$ws.Range("E7").Value = '-3.52%'
$ws.Range("D8").Value = '1.905'
$ws.Range("E8").Value = '-4.15%'
$ws.Range("D9").Value = '2.777'
$ws.Range("E9").Value = '-4.44%'
$ws.Range("D10").Value = '0.9441'
$ws.Range("E10").Value = '2.00%'
$ws.Range("D11").Value = '0.1175'
$ws.Range("E11").Value = '-8.49%'
$ws.Range("D12").Value = '0.1892'
$ws.Range("E12").Value = '-3.95%'
$ws.Range("D13").Value = '0.09729'
$ws.Range("E13").Value = '3.20%'
$ws.Range("D14").Value = '0.04305'
$ws.Range("E14").Value = '11.59%'
$ws.Range("D15").Value = '0.1069'
$ws.Range("E15").Value = '0.75%'
$ws.Range("E16").Value = '-1.95%'
$ws.Range("D17").Value = '0.006075'
$ws.Range("E17").Value = '-0.49%'
$ws.Range("D18").Value = '3.556'
$ws.Range("E18").Value = '3.33%'
$ws.Range("D19").Value = '0.3517'
$ws.Range("E19").Value = '-0.54%'
$ws.Range("D20").Value = '8.811'
$ws.Range("E20").Value = '-2.64%'
$ws.Range("D21").Value = '0.1359'
$ws.Range("E21").Value = '-0.25%'
$ws.Range("D22").Value = '0.2601'
$ws.Range("E22").Value = '5.81%'
$ws.Range("D23").Value = '0.04385'
$ws.Range("E23").Value = '-0.56%'
$ws.Range("E24").Value = '-2.49%'
$ws.Range("D25").Value = '0.004328'
$ws.Range("E25").Value = '-1.73%'
$ws.Range("D26").Value = '0.0001239'
$ws.Range("E26").Value = '3.24%'
$ws.Range("D27").Value = '0.0004018'
$ws.Range("E27").Value = '31.95%'
$ws.Range("D39").Value = '0.02650'
$ws.Range("E39").Value = '-6.37%'
$ws.Range("D40").Value = '0.05576'
$ws.Range("E40").Value = '1.24%'
$ws.Range("D41").Value = '0.007890'
$ws.Range("E41").Value = '1.07%'
$ws.Range("D42").Value = '0.009775'
$ws.Range("E42").Value = '5.02%'
$ws.Range("D43").Value = '0.1406'
$ws.Range("E43").Value = '-2.04%'
$ws.Range("D44").Value = '0.002127'
$ws.Range("E44").Value = '-0.62%'
$ws.Range("D45").Value = '0.009622'
$ws.Range("E45").Value = '-13.18%'
$ws.Range("D46").Value = '0.00007338'
$ws.Range("E46").Value = '3.56%'
$ws.Range("D47").Value = '0.00000000755'
$ws.Range("E47").Value = '0.73%'
$ws.Range("D48").Value = '0.003470'
$ws.Range("E48").Value = '2.06%'
$ws.Range("D49").Value = '0.002287'
$ws.Range("E49").Value = '0.38%'
$ws.Range("D50").Value = '0.00002115'
$ws.Range("E50").Value = '0.73%'
$ws.Range("D51").Value = '0.0002015'
$ws.Range("E51").Value = '0.73%'

# Restore the default (unstyled) look so no stray number-format is left behind
foreach ($addr in $forceCells) {
    $ws.Range($addr).Style = "Normal"
}
